# Weekly update for "Hortaliza, Terminal Hortofrutícola Agro Chillán - Cebollín".
# A new week of data lands at the top of this price series and every
# existing observation shifts down one row (row 122 keeps its own values;
# rows 123-184 each inherit the values that used to sit one row above them;
# a brand-new row 185 is appended carrying what used to be row 184's data).
#
# Columns A, B, C, E, F, G, H and R are constant for this whole
# market/category block, so only D (Fecha), I (Calidad), J:Q (Volumen ...
# Kg o Unidades) actually need to move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the date column down by one row.
$ws.Range("D123:D185").Value2 = $ws.Range("D122:D184").Value2

# Shift the quality column down by one row.
$ws.Range("I123:I185").Value2 = $ws.Range("I122:I184").Value2

# Shift the contiguous Volumen..Kg-o-Unidades block (J:Q) down by one row.
$ws.Range("J123:Q185").Value2 = $ws.Range("J122:Q184").Value2

# Row 185 is brand new - backfill the columns that are constant across the
# whole block (A, B, C, E, F, G, H, R) from the row right above it.
$ws.Range("A185").Value2 = $ws.Range("A184").Value2
$ws.Range("B185").Value2 = $ws.Range("B184").Value2
$ws.Range("C185").Value2 = $ws.Range("C184").Value2
$ws.Range("E185").Value2 = $ws.Range("E184").Value2
$ws.Range("F185").Value2 = $ws.Range("F184").Value2
$ws.Range("G185").Value2 = $ws.Range("G184").Value2
$ws.Range("H185").Value2 = $ws.Range("H184").Value2
$ws.Range("R185").Value2 = $ws.Range("R184").Value2

# D185 is a freshly-created cell so it doesn't inherit the date-time number
# format the rest of column D uses; copy it explicitly.
$ws.Range("D185").NumberFormat = $ws.Range("D184").NumberFormat
